$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Inventory
$ws.Range("C4").Value = 347000000
$ws.Range("D4").Value = 365000000
$ws.Range("E4").Value = 392000000
$ws.Range("F4").Value = 401000000

# Row 13 - Accounts Payable
$ws.Range("C13").Value = 229000000
$ws.Range("D13").Value = 217000000
$ws.Range("E13").Value = 245000000
$ws.Range("F13").Value = 263000000

# Row 25 - Additional Paid In Capital (B25 was an empty inline string, now numeric)
$ws.Range("B25").Value = 398400000

# Row 26 - Common Stock (Net)
$ws.Range("B26").Value = 400000

# Row 27 - Retained Earnings
$ws.Range("B27").Value = 3412600000

# Row 28 - Treasury Stock
$ws.Range("B28").Value = 50200000

# Row 32 - Shares (Common)
$ws.Range("B32").Value = 37061000

# Row 33 - Shareholders Equity (Tangible)
$ws.Range("B33").Value = 799099600

# Row 34 - Net Debt
$ws.Range("G34").Value = 651100000

# Row 35 - Total Debt
$ws.Range("G35").Value = 850600000
